{"js": "// Update the \"Overall pairwise comparisons\" table:\n//  - retitle it from the ITS/Phytohormone experiment to the 16S/Insect\n//    Herbivores experiment\n//  - replace the numeric results (SumsOfSqs, MeanSqs, F.Model, R2,\n//    P.value, P.value.corrected) for every data row with the new values\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\n\n// Helper: replace the text of the first paragraph/run found in a cell\n// while leaving paragraph/run formatting untouched.\nfunction setCellText(rowIndex, colIndex, newText) {\n  const cell = table.getCell(rowIndex, colIndex);\n  const para = cell.body.paragraphs.getFirst();\n  const rng = para.getRange();\n  rng.insertText(newText, Word.InsertLocation.replace);\n}\n\n// Row 0 (merged title row) / column 0 holds the table caption.\nsetCellText(\n  0,\n  0,\n  \"Table X. Overall pairwise comparisons between treatments of the 16S community for the Insect Herbivores experiment \"\n);\n\n// Data rows start at index 2 (index 0 = title, index 1 = header row).\n// Columns: 0 distance, 1 combination, 2 SumsOfSqs, 3 MeanSqs, 4 F.Model,\n// 5 R2, 6 P.value, 7 P.value.corrected.\nconst dataUpdates = [\n  // [rowIndex, colIndex, newValue]\n  [2, 2, \"0.08\"],\n  [2, 3, \"0.08\"],\n  [2, 4, \"0.91\"],\n  [2, 5, \".03\"],\n  [2, 6, \"0.68\"],\n  [2, 7, \"0.68\"],\n\n  [3, 2, \"0.10\"],\n  [3, 3, \"0.10\"],\n  [3, 4, \"1.22\"],\n  [3, 5, \".04\"],\n  [3, 6, \"0.11\"],\n  [3, 7, \"0.17\"],\n\n  [4, 2, \"0.10\"],\n  [4, 3, \"0.10\"],\n  [4, 4, \"1.24\"],\n  [4, 5, \".04\"],\n  [4, 6, \"0.10\"],\n  [4, 7, \"0.17\"],\n\n  [5, 2, \"0.10\"],\n  [5, 3, \"0.10\"],\n  [5, 4, \"1.05\"],\n  [5, 5, \".04\"],\n  [5, 6, \"0.33\"],\n  [5, 7, \"0.33\"],\n\n  [6, 2, \"0.14\"],\n  [6, 3, \"0.14\"],\n  [6, 4, \"1.58\"],\n  [6, 5, \".05\"],\n  [6, 6, \"0.00\"],\n\n  [7, 2, \"0.12\"],\n  [7, 3, \"0.12\"],\n  [7, 4, \"1.39\"],\n  [7, 5, \".05\"],\n  [7, 6, \"0.04\"],\n  [7, 7, \"0.06\"],\n];\n\nfor (const [rowIndex, colIndex, newValue] of dataUpdates) {\n  setCellText(rowIndex, colIndex, newValue);\n}\n\nawait context.sync();\n", "ps1": "# Update the \"Overall pairwise comparisons\" table:\n#  - retitle it from the ITS/Phytohormone experiment to the 16S/Insect\n#    Herbivores experiment\n#  - replace the numeric results (SumsOfSqs, MeanSqs, F.Model, R2,\n#    P.value, P.value.corrected) for every data row with the new values\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n# Row 1 is the merged title row spanning all 8 columns.\n$t.Cell(1, 1).Range.Text = \"Table X. Overall pairwise comparisons between treatments of the 16S community for the Insect Herbivores experiment \"\n\n# Data rows are table rows 3-8 (row 1 = title, row 2 = header).\n# Columns: 1 distance, 2 combination, 3 SumsOfSqs, 4 MeanSqs, 5 F.Model,\n# 6 R2, 7 P.value, 8 P.value.corrected.\n$updates = @(\n    @(3, 3, \"0.08\"), @(3, 4, \"0.08\"), @(3, 5, \"0.91\"), @(3, 6, \".03\"), @(3, 7, \"0.68\"), @(3, 8, \"0.68\"),\n    @(4, 3, \"0.10\"), @(4, 4, \"0.10\"), @(4, 5, \"1.22\"), @(4, 6, \".04\"), @(4, 7, \"0.11\"), @(4, 8, \"0.17\"),\n    @(5, 3, \"0.10\"), @(5, 4, \"0.10\"), @(5, 5, \"1.24\"), @(5, 6, \".04\"), @(5, 7, \"0.10\"), @(5, 8, \"0.17\"),\n    @(6, 3, \"0.10\"), @(6, 4, \"0.10\"), @(6, 5, \"1.05\"), @(6, 6, \".04\"), @(6, 7, \"0.33\"), @(6, 8, \"0.33\"),\n    @(7, 3, \"0.14\"), @(7, 4, \"0.14\"), @(7, 5, \"1.58\"), @(7, 6, \".05\"), @(7, 7, \"0.00\"),\n    @(8, 3, \"0.12\"), @(8, 4, \"0.12\"), @(8, 5, \"1.39\"), @(8, 6, \".05\"), @(8, 7, \"0.04\"), @(8, 8, \"0.06\")\n)\n\nforeach ($u in $updates) {\n    $row = $u[0]\n    $col = $u[1]\n    $val = $u[2]\n    $t.Cell($row, $col).Range.Text = $val\n}\n"}
